$d = $word.ActiveDocument

$replacements = @(
    @{old="69×48="; new="95×60="},
    @{old="71×87="; new="40×81="},
    @{old="29×98="; new="99×28="},
    @{old="66×73="; new="27×77="},
    @{old="68×68="; new="78×83="},
    @{old="81×47="; new="53×29="},
    @{old="54×39="; new="51×28="},
    @{old="17×85="; new="58×27="},
    @{old="25×51="; new="50×77="},
    @{old="73×20="; new="59×70="},
    @{old="71×34="; new="27×51="},
    @{old="27×34="; new="71×12="},
    @{old="12×23="; new="69×20="},
    @{old="93×59="; new="25×43="},
    @{old="86×43="; new="88×32="},
    @{old="25×74="; new="91×70="},
    @{old="89×87="; new="76×71="},
    @{old="75×35="; new="19×20="},
    @{old="34×27="; new="13×85="},
    @{old="19×25="; new="25×86="},
    @{old="13×67="; new="57×71="},
    @{old="93×53="; new="68×48="},
    @{old="65×57="; new="58×17="},
    @{old="69×91="; new="86×67="},
    @{old="59×32="; new="55×39="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
